$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Move column N's values (In Advance/Late header split + data) one column to the right,
# into the (previously empty) column O, leaving N blank - this mirrors an "Insert Cut Cells"
# style move of just that column's content rather than a full sheet-wide column insert.
$ws.Range("N1:N14").Cut($ws.Range("O1:O14"))

# Resize the columns: the (now blank) column N takes on a fresh width, the relocated
# column (now O) takes the old column-N width, and the following column (P) takes the
# old column-O width.
$ws.Columns("N").ColumnWidth = 10.333333333333334
$ws.Columns("O").ColumnWidth = 4.166666666666667
$ws.Columns("P").ColumnWidth = 11

# Make "Repayment Schedule" the active sheet and select cell M18 on it (this also clears
# the previous tab-selected/selection state on the "Transactions" sheet).
$ws.Activate()
$ws.Range("M18").Select()
